$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update shortname (B2) first so the new "255e" shared string is allocated
# before the new product-name string, matching the target shared-string order.
$ws1.Range("B2").Value = "255e"

# Update the product name on both sheets to the new value.
$ws1.Range("B1").Value = "2555-MS-EI-DB-SAR-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-PER-1st"
$ws2.Range("B1").Value = "2555-MS-EI-DB-SAR-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-PER-1st"

# Move the active selection on ProductLoanInput from B15 to B1.
$ws1.Range("B1").Select()

# Make ProductLoanOutput the active sheet/tab (selection stays at B1 there).
$ws2.Activate()
